$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.246.27'
$ws.Range("E2").Value = '  -1.76%  '
$ws.Range("D3").Value = '3.000.08'
$ws.Range("E3").Value = '  +0.09%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Formula = "'508.31"
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("D6").Formula = "'138.46"
$ws.Range("E6").Value = '  +1.23%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Formula = "'0.429"
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").Value = '  -2.32%  '
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("E11").Value = '  +3.29%  '
$ws.Range("D12").Value = '3.518.93'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("E14").Value = '  -2.24%  '
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").Value = '56.206.96'
$ws.Range("E16").Value = '  -1.88%  '
$ws.Range("D17").Value = '3.007.90'
$ws.Range("E17").Value = '  +0.31%  '
$ws.Range("D18").Formula = "'5.93"
$ws.Range("E18").Value = '  -2.29%  '
$ws.Range("D19").Formula = "'12.89"
$ws.Range("E19").Value = '  +2.30%  '
$ws.Range("D20").Formula = "'7.98"
$ws.Range("E20").Value = '  +1.77%  '
$ws.Range("D21").Formula = "'331.78"
$ws.Range("E21").Value = '  +3.78%  '
$ws.Range("D22").Formula = "'0.998"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("D23").Formula = "'0.496"
$ws.Range("E23").Value = '  +0.91%  '
$ws.Range("D24").Formula = "'64.68"
$ws.Range("E24").Value = '  +2.32%  '
$ws.Range("D25").Value = '3.131.92'
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("D26").Formula = "'0.165"
$ws.Range("E26").Value = '  +1.25%  '
$ws.Range("D27").Formula = "'1.00"
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = '0.0₃0934'
$ws.Range("E28").Value = '  +5.11%  '
$ws.Range("D29").Formula = "'6.33"
$ws.Range("E29").Value = '  -4.31%  '
$ws.Range("E30").Value = '  -2.95%  '
$ws.Range("E31").Value = '  +0.61%  '
$ws.Range("D32").Formula = "'20.29"
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("D34").Formula = "'152.46"
$ws.Range("E34").Value = '  -1.72%  '
$ws.Range("D35").Formula = "'4.42"
$ws.Range("E35").Value = '  -2.57%  '
$ws.Range("D36").Formula = "'26.61"
$ws.Range("E36").Value = '  +9.49%  '
$ws.Range("D37").Formula = "'5.80"
$ws.Range("E37").Value = '  +0.71%  '
$ws.Range("E38").Value = '  -0.68%  '
$ws.Range("D39").Formula = "'0.0659"
$ws.Range("E39").Value = '  -0.19%  '
$ws.Range("D40").Value = '3.039.64'
$ws.Range("E40").Value = '  +0.40%  '
$ws.Range("D41").Formula = "'36.39"
$ws.Range("E41").Value = '  -3.87%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").Formula = "'3.77"
$ws.Range("E43").Value = '  +1.33%  '
$ws.Range("D44").Formula = "'0.655"
$ws.Range("E44").Value = '  +1.47%  '
$ws.Range("D45").Value = '2.194.92'
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("E46").Value = '  -2.62%  '
$ws.Range("D47").Formula = "'0.0239"
$ws.Range("E47").Value = '  +2.45%  '
$ws.Range("D48").Formula = "'5.82"
$ws.Range("E48").Value = '  -1.83%  '
$ws.Range("D49").Formula = "'0.918"
$ws.Range("E49").Value = '  -1.64%  '
$ws.Range("D50").Formula = "'19.47"
$ws.Range("E50").Value = '  +1.42%  '
$ws.Range("D51").Formula = "'0.0847"
$ws.Range("E51").Value = '  -2.08%  '
